# Lab2.xlsx — "Refactoring and completing task 3"
#
# 1) task2: selection collapses from A1:D5 to a single active cell A6.
# 2) task3: the results table (B3:D5) gets filled in with the measured
#    timings for the three vector sizes (previously empty placeholder
#    cells), two helper columns (B, C) get explicit widths, and the
#    worksheet selection moves to K8.

$wb = $excel.ActiveWorkbook

# --- task2: just a selection/cursor move, no data changes -----------------
$ws2 = $wb.Worksheets.Item("task2")
$ws2.Activate()
[void]$ws2.Range("A6").Select()

# --- task3: fill in the completed measurements -----------------------------
$ws3 = $wb.Worksheets.Item("task3")
$ws3.Activate()

# Row 3 (size 512) — written left-to-right.
$ws3.Range("B3").Value = "0.0176926"
$ws3.Range("C3").Value = "1.28419"
$ws3.Range("D3").Value = "1.37147"

# Row 4 (size 1024) — written right-to-left so the new shared-string table
# keeps the same first-use ordering as the authored workbook.
$ws3.Range("D4").Value = "7.35719"
$ws3.Range("C4").Value = "7.73854"
$ws3.Range("B4").Value = "0.138595"

# Row 5 (size 2048) — written left-to-right.
$ws3.Range("B5").Value = "2.02386"
$ws3.Range("C5").Value = "31.7538"
$ws3.Range("D5").Value = "31.4294"

# New explicit widths for the two middle columns used by the table.
$ws3.Columns.Item(2).ColumnWidth = 8.6667
$ws3.Columns.Item(3).ColumnWidth = 10

# Leave the cursor on K8, as in the saved file.
[void]$ws3.Range("K8").Select()
